# Update countries & provincias Spain
# Japon's case count overtook Chile's, and Argentina's overtook Sudafrica's,
# so those pairs of rows swap country labels (values follow their country).
# A handful of other rows (Noruega, etc.) simply get refreshed numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Row($row, $name, $values) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = $values[0]
    $ws.Cells.Item($row, 3).Value = $values[1]
    $ws.Cells.Item($row, 4).Value = $values[2]
    $ws.Cells.Item($row, 5).Value = $values[3]
    $ws.Cells.Item($row, 6).Value = $values[4]
    $ws.Cells.Item($row, 7).Value = $values[5]
    $ws.Cells.Item($row, 8).Value = $values[6]
}

# Row 27 now holds Japon's updated figures (Japon passed Chile in total cases)
Set-Row 27 "Japon" @(9231, 605, 935, 8106, 193, 12, 190)

# Row 28 now holds Chile's figures (unchanged totals, just shifted down a row)
Set-Row 28 "Chile" @(8807, 534, 3299, 5403, 384, 11, 105)

# Row 33 (Noruega) refreshed figures, no reordering involved
Set-Row 33 "Noruega" @(6905, 108, 32, 6721, 64, 2, 152)

# Row 55 now holds Argentina's updated figures (Argentina passed Sudafrica)
Set-Row 55 "Argentina" @(2669, 98, 631, 1916, 121, 10, 122)

# Row 56 now holds Sudafrica's figures (unchanged totals, just shifted down a row)
Set-Row 56 "Sudafrica" @(2605, 99, 903, 1654, 7, 14, 48)
